$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.603.93"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.444.60"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "547.56"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "145.73"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").Value = "  +0.04%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").Value = "2.443.55"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  -2.49%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.86"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "2.882.61"
$ws.Range("E15").Value = "  -2.36%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000167"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "61.841.50"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "2.438.30"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("E19").Value = "  -4.09%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  -2.13%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "318.21"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E23").Value = "  -0.13%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.99%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "63.90"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "0.0₃0970"
$ws.Range("E26").Value = "  -6.71%  "
$ws.Range("D27").Value = "2.564.60"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  +1.40%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "523.85"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.63"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.27%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("E39").Value = "  +0.52%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.19"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("E41").Value = "  +1.07%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "138.53"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("E43").Value = "  +0.18%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "40.36"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -4.06%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "142.46"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.73%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -0.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0524"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("E50").Value = "  -0.85%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.82%  "
